$wb = $excel.ActiveWorkbook

$wsEntities = $wb.Worksheets.Item("Entities")
$wsMethods  = $wb.Worksheets.Item("Methods")

# Mark the QuestionTimeline entity row as "Done" (apply the built-in "Good" style)
$wsEntities.Range("A12").Style = "Good"

# Mark each already-documented method row as "Done" (built-in "Good" style)
$doneRows = @(6,7,8,9,10,11,21,22,23,24,28,29,30,31,32,33,34,35,36)
foreach ($r in $doneRows) {
    $cell = $wsMethods.Range("A$r")
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# Update selections on each sheet
$wsEntities.Range("A12").Select()
$wsMethods.Range("A36").Select()

# Methods becomes the active/selected tab
$wsMethods.Activate()
